$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The task list previously ending at row 92 ("...significant items more
# clearly #64") is bumped up one row, and the now-vacated slot plus the
# newly-added tasks below get their hours filled in / a new task added.

$ws.Range("D91").Value = "As a user, I want to see the significant items more clearly#64"
$ws.Range("D92").Value = ""

$ws.Range("E93").Value = 1
$ws.Range("G93").Value = 8

$ws.Range("E94").Value = 0.5

$ws.Range("E95").Value = 2

$ws.Range("E96").Value = 2

$ws.Range("E97").Value = 0.5

$nbsp = [char]0x00A0
$ws.Range("D98").Value = "As a coder, I want to fix the cursor not appearing on screen correctly${nbsp}#68"
$ws.Range("E98").Value = 2

$ws.Range("D99").Value = ""

$ws.Range("E100").Value = 1
$ws.Range("G100").Value = 7

# Reflect where the sheet was scrolled / selected when the edits were made.
$ws.Range("H111").Select()
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 2
